$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2134.6667
$ws.Range("I15").Value = 2134.6667
$ws.Range("K15").Value = 6404.000100000001
$ws.Range("M15").Value = -6235.000100000001
$ws.Range("H33").Value = 6565.125
$ws.Range("I33").Value = 10187.2
$ws.Range("J33").Value = 528.3333
$ws.Range("K33").Value = 10187.2
$ws.Range("L33").Value = 528.3333
$ws.Range("M33").Value = -9958.200000000001
$ws.Range("N33").Value = -986.3333
$ws.Range("H40").Value = 4394.7646
$ws.Range("J40").Value = 4622.643
$ws.Range("L40").Value = 4622.643
$ws.Range("N40").Value = -4972.643
$ws.Range("H96").Value = 4273.8
$ws.Range("I96").Value = 860
$ws.Range("J96").Value = 6549.6665
$ws.Range("K96").Value = 2580
$ws.Range("L96").Value = 19648.9995
$ws.Range("M96").Value = -1207
$ws.Range("N96").Value = -22394.9995
$ws.Range("H103").Value = 1200
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -5672
$ws.Range("H106").Value = 12061.333
$ws.Range("I106").Value = 8473.6
$ws.Range("K106").Value = 8473.6
$ws.Range("M106").Value = -7842.6
$ws.Range("H112").Value = 2335.6667
$ws.Range("I112").Value = 2625
$ws.Range("J112").Value = 2299.5
$ws.Range("K112").Value = 7875
$ws.Range("L112").Value = 6898.5
$ws.Range("M112").Value = -6767
$ws.Range("N112").Value = -9114.5
$ws.Range("H116").Value = 7624.4165
$ws.Range("J116").Value = 7812.125
$ws.Range("L116").Value = 7812.125
$ws.Range("N116").Value = -14696.125
$ws.Range("H129").Value = 799.875
$ws.Range("I129").Value = 799.875
$ws.Range("K129").Value = 2399.625
$ws.Range("M129").Value = 2600.375
$ws.Range("H132").Value = 5307.628
$ws.Range("I132").Value = 4057.718
$ws.Range("K132").Value = 12173.154
$ws.Range("M132").Value = -9643.153999999999
$ws.Range("H138").Value = 4810.516
$ws.Range("J138").Value = 5216.5317
$ws.Range("L138").Value = 15649.5951
$ws.Range("N138").Value = -25929.5951

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41018.227
$ws.Range("I32").Value = 45656.348
$ws.Range("J32").Value = 16900
$ws.Range("K32").Value = 45656.348
$ws.Range("L32").Value = 16900
$ws.Range("M32").Value = -45369.348
$ws.Range("N32").Value = -17474
$ws.Range("H61").Value = 2974.4119
$ws.Range("I61").Value = 3111.9285
$ws.Range("J61").Value = 2332.6667
$ws.Range("K61").Value = 3111.9285
$ws.Range("L61").Value = 2332.6667
$ws.Range("M61").Value = -2899.9285
$ws.Range("N61").Value = -2756.6667
$ws.Range("H88").Value = 1703.4375
$ws.Range("I88").Value = 1086.125
$ws.Range("K88").Value = 1086.125
$ws.Range("M88").Value = -680.125
$ws.Range("H91").Value = 1703.4375
$ws.Range("I91").Value = 1086.125
$ws.Range("K91").Value = 1086.125
$ws.Range("M91").Value = 317.875
$ws.Range("H136").Value = 2974.4119
$ws.Range("I136").Value = 3111.9285
$ws.Range("J136").Value = 2332.6667
$ws.Range("K136").Value = 9335.7855
$ws.Range("L136").Value = 6998.000100000001
$ws.Range("M136").Value = -6785.7855
$ws.Range("N136").Value = -12098.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3832.9678
$ws.Range("I105").Value = 3718.0688
$ws.Range("K105").Value = 3718.0688
$ws.Range("M105").Value = -1971.0688
$ws.Range("H134").Value = 1602.1072
$ws.Range("I134").Value = 1637
$ws.Range("J134").Value = 660
$ws.Range("K134").Value = 4911
$ws.Range("L134").Value = 1980
$ws.Range("M134").Value = -2376
$ws.Range("N134").Value = -7050

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H105").Value = 652.7143
$ws.Range("I105").Value = 652.7143
$ws.Range("K105").Value = 652.7143
$ws.Range("M105").Value = 1094.2857
$ws.Range("H122").Value = 1555.8695
$ws.Range("I122").Value = 1558.8125
$ws.Range("J122").Value = 1549.1428
$ws.Range("K122").Value = 4676.4375
$ws.Range("L122").Value = 4647.428400000001
$ws.Range("M122").Value = -2226.4375
$ws.Range("N122").Value = -9547.428400000001
$ws.Range("H132").Value = 1719.3077
$ws.Range("I132").Value = 1546.5278
$ws.Range("J132").Value = 3792.6667
$ws.Range("K132").Value = 4639.5834
$ws.Range("L132").Value = 11378.0001
$ws.Range("M132").Value = -2109.5834
$ws.Range("N132").Value = -16438.0001
$ws.Range("H134").Value = 45404.695
$ws.Range("I134").Value = 48776.57
$ws.Range("K134").Value = 146329.71
$ws.Range("M134").Value = -143794.71

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.53846
$ws.Range("I2").Value = 38
$ws.Range("K2").Value = 228
$ws.Range("M2").Value = -115
$ws.Range("H23").Value = 425.5
$ws.Range("I23").Value = 70.5
$ws.Range("K23").Value = 211.5
$ws.Range("M23").Value = 23.5
$ws.Range("H36").Value = 498
$ws.Range("I36").Value = 498
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1494
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1325
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 330.5
$ws.Range("I38").Value = 390
$ws.Range("K38").Value = 1170
$ws.Range("M38").Value = -823
$ws.Range("H62").Value = 7999.5
$ws.Range("J62").Value = 7999.5
$ws.Range("L62").Value = 23998.5
$ws.Range("N62").Value = -25370.5
$ws.Range("H65").Value = 7999.5
$ws.Range("J65").Value = 7999.5
$ws.Range("L65").Value = 71995.5
$ws.Range("N65").Value = -78859.5
$ws.Range("H97").Value = 353.45456
$ws.Range("I97").Value = 273.25
$ws.Range("J97").Value = 399.2857
$ws.Range("K97").Value = 819.75
$ws.Range("L97").Value = 1197.8571
$ws.Range("M97").Value = -323.75
$ws.Range("N97").Value = -2189.8571

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 20001
$ws.Range("K70").Value = 20001
$ws.Range("M70").Value = -19731
$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 20001
$ws.Range("K73").Value = 20001
$ws.Range("M73").Value = -19065
$ws.Range("H80").Value = 2777.5334
$ws.Range("J80").Value = 2918.7778
$ws.Range("L80").Value = 2918.7778
$ws.Range("N80").Value = -4914.7778
$ws.Range("H83").Value = 2777.5334
$ws.Range("J83").Value = 2918.7778
$ws.Range("L83").Value = 14593.889
$ws.Range("N83").Value = -24577.889
$ws.Range("H102").Value = 2685.348
$ws.Range("I102").Value = 1888.2
$ws.Range("K102").Value = 1888.2
$ws.Range("M102").Value = -266.2
$ws.Range("H122").Value = 2941.2083
$ws.Range("I122").Value = 2959.7
$ws.Range("J122").Value = 2848.75
$ws.Range("K122").Value = 8879.099999999999
$ws.Range("L122").Value = 8546.25
$ws.Range("M122").Value = -6429.099999999999
$ws.Range("N122").Value = -13446.25
$ws.Range("H132").Value = 66594.625
$ws.Range("I132").Value = 73965.36
$ws.Range("J132").Value = 14999.5
$ws.Range("K132").Value = 221896.08
$ws.Range("L132").Value = 44998.5
$ws.Range("M132").Value = -219366.08
$ws.Range("N132").Value = -50058.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 741.63635
$ws.Range("I16").Value = 556.55554
$ws.Range("J16").Value = 1574.5
$ws.Range("K16").Value = 556.55554
$ws.Range("L16").Value = 1574.5
$ws.Range("M16").Value = -386.55554
$ws.Range("N16").Value = -1914.5
$ws.Range("H22").Value = 52728.3
$ws.Range("I22").Value = 143714
$ws.Range("K22").Value = 143714
$ws.Range("M22").Value = -143419
$ws.Range("H27").Value = 52728.3
$ws.Range("I27").Value = 143714
$ws.Range("K27").Value = 143714
$ws.Range("M27").Value = -143607
$ws.Range("H40").Value = 4633.25
$ws.Range("I40").Value = 3444.3333
$ws.Range("J40").Value = 8200
$ws.Range("K40").Value = 3444.3333
$ws.Range("L40").Value = 8200
$ws.Range("M40").Value = -3308.3333
$ws.Range("N40").Value = -8472
$ws.Range("H133").Value = 67997
$ws.Range("J133").Value = 67997
$ws.Range("L133").Value = 67997
$ws.Range("N133").Value = -73057

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2070
$ws.Range("I81").Value = 790.1667
$ws.Range("J81").Value = 3989.75
$ws.Range("K81").Value = 1580.3334
$ws.Range("L81").Value = 7979.5
$ws.Range("M81").Value = -519.3334
$ws.Range("N81").Value = -10101.5
$ws.Range("H84").Value = 2070
$ws.Range("I84").Value = 790.1667
$ws.Range("J84").Value = 3989.75
$ws.Range("K84").Value = 7901.666999999999
$ws.Range("L84").Value = 39897.5
$ws.Range("M84").Value = -2597.666999999999
$ws.Range("N84").Value = -50505.5
$ws.Range("H107").Value = 608.7273
$ws.Range("I107").Value = 608.7273
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1826.1819
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 93.81809999999996
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 725.6667
$ws.Range("I122").Value = 764.0909
$ws.Range("K122").Value = 2292.2727
$ws.Range("M122").Value = 157.7273
$ws.Range("H136").Value = 3777.1628
$ws.Range("I136").Value = 3918.8438
$ws.Range("J136").Value = 3365
$ws.Range("K136").Value = 11756.5314
$ws.Range("L136").Value = 10095
$ws.Range("M136").Value = -9206.5314
$ws.Range("N136").Value = -15195
